$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - match formatting of existing header row (bold font, thin box border, centered/top aligned)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

foreach ($addr in @("I1", "J1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4160    # xlTop

    $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft, xlContinuous
    $cell.Borders.Item(7).Weight = 2      # xlThin
    $cell.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $cell.Borders.Item(8).Weight = 2
    $cell.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $cell.Borders.Item(9).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $cell.Borders.Item(10).Weight = 2
}

# Data values for I2:J24
$values = @{
    2  = @(8, 8)
    3  = @(9, 9)
    4  = @(8, 8)
    5  = @(10, 10)
    6  = @(9, 9)
    7  = @(6, 7)
    8  = @(5, 5)
    9  = @(8, 9)
    10 = @(7, 7)
    11 = @(6, 6)
    12 = @(7, 7)
    13 = @(9, 9)
    14 = @(6, 6)
    15 = @(8, 8)
    16 = @(7, 7)
    17 = @(6, 6)
    18 = @(5, 5)
    19 = @(8, 8)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(4, 4)
    23 = @(4, 4)
    24 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
